$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 70
$prev = 69

# Column A ("Indice") reuses the bold/bordered header-row style from the
# row above it - copy the source cell (format + value) then overwrite
# the value, which keeps the existing shared style instead of minting one.
$ws.Cells.Item($prev, 1).Copy($ws.Cells.Item($r, 1))
$ws.Cells.Item($r, 1).Value = 69

$ws.Cells.Item($r, 2).Value = "morocco"
$ws.Cells.Item($r, 3).Value = "botola-pro"
$ws.Cells.Item($r, 4).Value = "2023-2024"

$ws.Cells.Item($r, 5).Value = 45254.83333333334
$ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($r, 6).Value = "Union Touarga"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = "FUS Rabat"
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 2.86
$ws.Cells.Item($r, 11).Value = "23/11/2023 02:12"
$ws.Cells.Item($r, 12).Value = 4.28
$ws.Cells.Item($r, 13).Value = "24/11/2023 19:50"
$ws.Cells.Item($r, 14).Value = 2.81
$ws.Cells.Item($r, 15).Value = "23/11/2023 02:12"
$ws.Cells.Item($r, 16).Value = 3.07
$ws.Cells.Item($r, 17).Value = "24/11/2023 19:50"
$ws.Cells.Item($r, 18).Value = 2.44
$ws.Cells.Item($r, 19).Value = "23/11/2023 02:12"
$ws.Cells.Item($r, 20).Value = 1.95
$ws.Cells.Item($r, 21).Value = "24/11/2023 19:50"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/morocco/botola-pro/union-touarga-fus-rabat/jDVARjxE/"
